$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5027842378199027
$ws.Cells.Item(2, 3).Value = 0.9899884981017091
$ws.Cells.Item(2, 4).Value = 0.5795766793119063
$ws.Cells.Item(2, 7).Value = 1.558460351833249
$ws.Cells.Item(2, 8).Value = 0.9990000000000001
$ws.Cells.Item(3, 2).Value = 0.2367633079409734
$ws.Cells.Item(3, 3).Value = 0.9953724551921432
$ws.Cells.Item(3, 4).Value = 0.3880481933457188
$ws.Cells.Item(3, 7).Value = 1.558460351833249
$ws.Cells.Item(3, 8).Value = 0.9990000000000001
$ws.Cells.Item(4, 2).Value = 0.279539502636717
$ws.Cells.Item(4, 3).Value = 0.994622281564909
$ws.Cells.Item(4, 4).Value = 0.4308435359205416
$ws.Cells.Item(4, 7).Value = 1.558460351833249
$ws.Cells.Item(4, 8).Value = 0.9990000000000001
$ws.Cells.Item(5, 2).Value = 0.4110342851853201
$ws.Cells.Item(5, 3).Value = 0.991894966652257
$ws.Cells.Item(5, 4).Value = 0.4923299625507078
$ws.Cells.Item(5, 7).Value = 1.558460351833249
$ws.Cells.Item(5, 8).Value = 0.9990000000000001
$ws.Cells.Item(6, 2).Value = 0.4921451733549996
$ws.Cells.Item(6, 3).Value = 0.9855444070527011
$ws.Cells.Item(6, 4).Value = 0.5158798569838562
$ws.Cells.Item(6, 7).Value = 1.558460351833249
$ws.Cells.Item(6, 8).Value = 0.9990000000000001
$ws.Cells.Item(7, 2).Value = 0.09825746411803755
$ws.Cells.Item(7, 3).Value = 0.9986420552412684
$ws.Cells.Item(7, 4).Value = 0.2566296982374509
$ws.Cells.Item(7, 7).Value = 1.558460351833249
$ws.Cells.Item(7, 8).Value = 0.9990000000000001
$ws.Cells.Item(8, 2).Value = 0.03278388427532332
$ws.Cells.Item(8, 3).Value = 0.9996621034951678
$ws.Cells.Item(8, 4).Value = 0.1307472531425428
$ws.Cells.Item(8, 7).Value = 1.558460351833249
$ws.Cells.Item(8, 8).Value = 0.9990000000000001
$ws.Cells.Item(9, 2).Value = 0.1107637560280684
$ws.Cells.Item(9, 3).Value = 0.9993402119087579
$ws.Cells.Item(9, 4).Value = 0.2294486971184508
$ws.Cells.Item(9, 7).Value = 1.558460351833249
$ws.Cells.Item(9, 8).Value = 0.9990000000000001
$ws.Cells.Item(10, 2).Value = 0.06723949699425752
$ws.Cells.Item(10, 3).Value = 0.9987716936268701
$ws.Cells.Item(10, 4).Value = 0.2052751533076876
$ws.Cells.Item(10, 7).Value = 1.558460351833249
$ws.Cells.Item(10, 8).Value = 0.9990000000000001
$ws.Cells.Item(11, 2).Value = 0.1311181845461133
$ws.Cells.Item(11, 3).Value = 0.9903125108657441
$ws.Cells.Item(11, 4).Value = 0.2835524564762872
$ws.Cells.Item(11, 7).Value = 1.558460351833249
$ws.Cells.Item(11, 8).Value = 0.9990000000000001
$ws.Cells.Item(12, 2).Value = 0.05224580291309232
$ws.Cells.Item(12, 3).Value = 0.9984635400702141
$ws.Cells.Item(12, 4).Value = 0.1677012470627132
$ws.Cells.Item(12, 7).Value = 1.558460351833249
$ws.Cells.Item(12, 8).Value = 0.9990000000000001
$ws.Cells.Item(13, 2).Value = 0.06020935002615127
$ws.Cells.Item(13, 3).Value = 0.9994288873323021
$ws.Cells.Item(13, 4).Value = 0.1787292105558562
$ws.Cells.Item(13, 7).Value = 1.558460351833249
$ws.Cells.Item(13, 8).Value = 0.9990000000000001
$ws.Cells.Item(14, 2).Value = 0.05551525225741943
$ws.Cells.Item(14, 3).Value = 0.9992483577529127
$ws.Cells.Item(14, 4).Value = 0.1886676226951151
$ws.Cells.Item(14, 7).Value = 1.558460351833249
$ws.Cells.Item(14, 8).Value = 0.9990000000000001
